$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

# Column -> new value mapping (applies identically to rows 2 and 3)
$updates = @{
    "D"  = -0.375
    "G"  = -2.138297872340425
    "H"  = -2.138297872340425
    "I"  = -2.148936170212766
    "J"  = -2.148936170212766
    "K"  = -12
    "L"  = -6.382978723404256
    "U"  = 2.46
    "V"  = 0.06406250000000001
    "W"  = -1.213346814964611
    "X"  = 0.08565043905558375
    "Y"  = -1.298997254020194
    "Z"  = 0.09606540623403168
    "AA" = -0.2064384261624936
    "AB" = 0.07627217525767803
    "AC" = -0.2827106014201716
    "AD" = 14.9
    "AF" = 14.9
    "AG" = 12.44
    "AH" = 0.2795497185741088
    "AI" = 0.2838095238095238
    "AJ" = 0.2446892210857592
    "AK" = 0.2486011191047162
    "AL" = 4.46
    "AM" = 4.429
    "AN" = -3.725
    "AO" = -0.9058295964125561
    "AP" = -3.11
    "AQ" = -0.9121697900203206
}

foreach ($row in 2..3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
